# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns and updates the Status text for the two localized-file rows on both the
# zh-cn and de-de worksheets, now that the handback has been completed.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- zh-cn sheet -----------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2: 2b431392-f2f5-4ead-aca5-bdecef3210cb.md
$wsZh.Range("B2").Value = $statusHandedBack
$wsZh.Range("E2").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.md"
$wsZh.Range("E2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/5602b44a9c5edf3c898dcb3c488843585b3b815e/e2e/2b431392-f2f5-4ead-aca5-bdecef3210cb.md", $null, $null, "2b431392-f2f5-4ead-aca5-bdecef3210cb.md") | Out-Null

$wsZh.Range("F2").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.zh-cn.xlf"
$wsZh.Range("F2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9a7a0dd33a5adadee11b78a9f8999ebe40ee760/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.zh-cn.xlf", $null, $null, "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.zh-cn.xlf") | Out-Null

$wsZh.Range("G2").Value = "2016-01-28 06:00:34"

# Row 3: 310c7c43-7abb-4f88-8dac-82b8c4bfe446.md
$wsZh.Range("B3").Value = $statusHandedBack
$wsZh.Range("E3").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md"
$wsZh.Range("E3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/5602b44a9c5edf3c898dcb3c488843585b3b815e/e2e/310c7c43-7abb-4f88-8dac-82b8c4bfe446.md", $null, $null, "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md") | Out-Null

$wsZh.Range("F3").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.zh-cn.xlf"
$wsZh.Range("F3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9a7a0dd33a5adadee11b78a9f8999ebe40ee760/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.zh-cn.xlf", $null, $null, "310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.zh-cn.xlf") | Out-Null

$wsZh.Range("G3").Value = "2016-01-28 06:00:34"

# --- de-de sheet -------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2: 2b431392-f2f5-4ead-aca5-bdecef3210cb.md
$wsDe.Range("B2").Value = $statusHandedBack
$wsDe.Range("E2").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.md"
$wsDe.Range("E2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/5602b44a9c5edf3c898dcb3c488843585b3b815e/e2e/2b431392-f2f5-4ead-aca5-bdecef3210cb.md", $null, $null, "2b431392-f2f5-4ead-aca5-bdecef3210cb.md") | Out-Null

$wsDe.Range("F2").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.de-de.xlf"
$wsDe.Range("F2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f101a9441c8ac7246ab2ca36a938c454b0364986/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.de-de.xlf", $null, $null, "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.de-de.xlf") | Out-Null

$wsDe.Range("G2").Value = "2016-01-28 06:00:53"

# Row 3: 310c7c43-7abb-4f88-8dac-82b8c4bfe446.md
$wsDe.Range("B3").Value = $statusHandedBack
$wsDe.Range("E3").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md"
$wsDe.Range("E3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/5602b44a9c5edf3c898dcb3c488843585b3b815e/e2e/310c7c43-7abb-4f88-8dac-82b8c4bfe446.md", $null, $null, "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md") | Out-Null

$wsDe.Range("F3").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.de-de.xlf"
$wsDe.Range("F3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f101a9441c8ac7246ab2ca36a938c454b0364986/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.de-de.xlf", $null, $null, "310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.de-de.xlf") | Out-Null

$wsDe.Range("G3").Value = "2016-01-28 06:00:53"

Write-Host "Handback report generated."
